$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $value) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws 'D2' '42.397.73'
$ws.Range('E2').Value = '  +0.63%  '
Set-TextValue $ws 'D3' '2.241.04'
$ws.Range('E3').Value = '  -0.14%  '
$ws.Range('E4').Value = '  +0.18%  '
Set-TextValue $ws 'D5' '245.41'
$ws.Range('E5').Value = '  -0.55%  '
Set-TextValue $ws 'D6' '0.629'
$ws.Range('E6').Value = '  +0.94%  '
Set-TextValue $ws 'D7' '74.96'
$ws.Range('E7').Value = '  -2.30%  '
$ws.Range('E8').Value = '  +0.16%  '
Set-TextValue $ws 'D9' '0.621'
$ws.Range('E9').Value = '  +0.95%  '
Set-TextValue $ws 'D10' '43.54'
$ws.Range('E10').Value = '  +3.79%  '
Set-TextValue $ws 'D11' '0.0958'
$ws.Range('E11').Value = '  +0.86%  '
Set-TextValue $ws 'D12' '7.15'
$ws.Range('E12').Value = '  +0.59%  '
$ws.Range('E13').Value = '  +0.24%  '
Set-TextValue $ws 'D14' '14.52'
$ws.Range('E14').Value = '  -1.94%  '
Set-TextValue $ws 'D15' '0.856'
$ws.Range('E15').Value = '  -0.58%  '
Set-TextValue $ws 'D16' '2.241.25'
$ws.Range('E16').Value = '  -1.18%  '
Set-TextValue $ws 'D17' '42.315.12'
$ws.Range('E17').Value = '  +0.78%  '
Set-TextValue $ws 'D18' '0.0000109'
$ws.Range('E18').Value = '  +11.15%  '
Set-TextValue $ws 'D19' '6.18'
$ws.Range('E19').Value = '  +1.00%  '
Set-TextValue $ws 'D20' '72.06'
$ws.Range('E20').Value = '  +0.29%  '
Set-TextValue $ws 'D21' '10.41'
$ws.Range('E21').Value = '  +42.79%  '
Set-TextValue $ws 'D22' '231.61'
$ws.Range('E22').Value = '  +0.62%  '
Set-TextValue $ws 'D23' '2.18'
$ws.Range('E23').Value = '  -6.57%  '
Set-TextValue $ws 'D24' '11.81'
$ws.Range('E24').Value = '  +3.85%  '
$ws.Range('E25').Value = '  +0.08%  '
Set-TextValue $ws 'D26' '3.65'
$ws.Range('E26').Value = '  +0.33%  '
Set-TextValue $ws 'D27' '2.31'
$ws.Range('E27').Value = '  +0.99%  '
$ws.Range('E28').Value = '  +6.73%  '
Set-TextValue $ws 'D29' '166.70'
$ws.Range('E29').Value = '  -1.98%  '
$ws.Range('E30').Value = '  +1.69%  '
Set-TextValue $ws 'D31' '5.90'
$ws.Range('E31').Value = '  +20.44%  '
Set-TextValue $ws 'D32' '0.0813'
$ws.Range('E32').Value = '  -1.52%  '
$ws.Range('E33').Value = '  -2.13%  '
$ws.Range('B34').Value = 'InjectiveProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws 'D34' '29.74'
$ws.Range('E34').Value = '  -10.65%  '
$ws.Range('B35').Value = 'Stellar'
$ws.Range('C35').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws 'D35' '0.125'
$ws.Range('E35').Value = '  +0.14%  '
Set-TextValue $ws 'D36' '4.56'
$ws.Range('E36').Value = '  +1.18%  '
Set-TextValue $ws 'D37' '0.0311'
$ws.Range('E37').Value = '  +2.91%  '
Set-TextValue $ws 'D38' '13.43'
$ws.Range('E38').Value = '  -6.31%  '
Set-TextValue $ws 'D39' '2.17'
$ws.Range('E39').Value = '  -0.72%  '
Set-TextValue $ws 'D40' '5.68'
$ws.Range('E40').Value = '  -3.64%  '
Set-TextValue $ws 'D41' '63.85'
$ws.Range('E41').Value = '  +4.66%  '
Set-TextValue $ws 'D42' '0.202'
$ws.Range('E42').Value = '  -0.76%  '
Set-TextValue $ws 'D43' '8.85'
$ws.Range('E43').Value = '  +2.05%  '
Set-TextValue $ws 'D44' '106.10'
$ws.Range('E44').Value = '  -6.28%  '
$ws.Range('E45').Value = '  +2.16%  '
Set-TextValue $ws 'D46' '0.996'
$ws.Range('E46').Value = '  -0.08%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws 'D47' '2.40'
$ws.Range('E47').Value = '  +4.22%  '
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws 'D48' '1.14'
$ws.Range('E48').Value = '  +0.44%  '
Set-TextValue $ws 'D49' '1.18'
$ws.Range('E49').Value = '  +1.13%  '
$ws.Range('E50').Value = '  +1.52%  '
Set-TextValue $ws 'D51' '4.13'
$ws.Range('E51').Value = '  -1.89%  '
